$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing table (A1:B7) down by one row, freeing row 1 for the
# new "section" labels that sit above each of the two side-by-side tables.
$ws.Rows("1:1").Insert()

# --- Section labels (row 1) ---
# Set E1 before A1 so the shared-string table ends up in the same order
# as the target workbook (All issuers, then Top performers only).
$ws.Range("E1").Value = "All issuers"
$ws.Range("A1").Value = "Top performers only"

# --- Column headers (row 2) ---
$ws.Range("C2").Value = "rel %"
$ws.Range("E2").Value = "violation_category"
$ws.Range("F2").Value = "count"
$ws.Range("G2").Value = "rel %"

# --- "All issuers" data table (columns E:G, rows 3-9) ---
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 166861
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 896413
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 229500
$ws.Range("E6").Value = 4
$ws.Range("F6").Value = 221699
$ws.Range("E7").Value = 5
$ws.Range("F7").Value = 535962
$ws.Range("E8").Value = 6
$ws.Range("F8").Value = 9273
$ws.Range("F9").Value = 331

# --- Relative percentage formulas ---
# "Top performers only" table (column C)
$ws.Range("C3:C8").NumberFormat = "0.00%"
$ws.Range("C3").Formula = '=B3/SUM($B$3:$B$8)'
$ws.Range("C4:C8").Formula = '=B4/SUM($B$3:$B$8)'

# "All issuers" table (column G)
$ws.Range("G3:G9").NumberFormat = "0.00%"
$ws.Range("G3").Formula = '=F3/SUM($F$3:$F$9)'
$ws.Range("G4:G9").Formula = '=F4/SUM($F$3:$F$9)'

# --- Misc sheet/workbook presentation tweaks ---
$ws.PageSetup.Orientation = 1
$ws.Range("H23").Select()
